$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column H ("Label")
$ws.Cells.Item(1, 8).Value = "Label"
$ws.Cells.Item(1, 8).Font.Bold = $true
$ws.Cells.Item(1, 8).HorizontalAlignment = -4108
$ws.Cells.Item(1, 8).VerticalAlignment = -4160
$ws.Cells.Item(1, 8).Borders.LineStyle = 1

# Populate Label column: 0 for Control rows, 1 for MDD rows
for ($r = 2; $r -le 21; $r++) {
    $diag = $ws.Cells.Item($r, 1).Value()
    if ($diag -like "Control*") {
        $ws.Cells.Item($r, 8).Value = 0
    } elseif ($diag -like "MDD*") {
        $ws.Cells.Item($r, 8).Value = 1
    }
}
